# Add two new "target_type" enum entries (ally_hp_highest, ally_hp_lowest)
# and rework the existing "ally" row into a "self" row, per the commit
# "[taekwon] enum 내 타겟 타입 추가" (add target type to enum).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 60: target_type / ally -> target_type / self -------------------
$ws.Cells.Item(60, 2).Value = "self"

# --- Row 61: target_type / enemy (unchanged content, kept as-is) --------
$ws.Cells.Item(61, 1).Value = "target_type"
$ws.Cells.Item(61, 2).Value = "enemy"
$ws.Cells.Item(61, 4).Value = 1

# --- New rows 62 & 63: copy formatting from row 61 first -----------------
$ws.Range("A61:E61").Copy()
$ws.Range("A62:E63").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the two new rows' group/enum names
$ws.Cells.Item(62, 1).Value = "target_type"
$ws.Cells.Item(63, 1).Value = "target_type"

# B column values: row 63 (ally_hp_lowest) entered before row 62 (ally_hp_highest)
$ws.Cells.Item(63, 2).Value = "ally_hp_lowest"
$ws.Cells.Item(62, 2).Value = "ally_hp_highest"

$ws.Cells.Item(62, 4).Value = 2
$ws.Cells.Item(63, 4).Value = 3

# --- Formulas (UPPER(group)&"_"&UPPER(enum)) for all four data rows -----
$ws.Cells.Item(60, 3).Formula = '=UPPER(A60)&"_"&UPPER(B60)'
$ws.Cells.Item(61, 3).Formula = '=UPPER(A61)&"_"&UPPER(B61)'
$ws.Cells.Item(62, 3).Formula = '=UPPER(A62)&"_"&UPPER(B62)'
$ws.Cells.Item(63, 3).Formula = '=UPPER(A63)&"_"&UPPER(B63)'

# --- Column E descriptions, filled in row order --------------------------
$ws.Cells.Item(60, 5).Value = "나 자신"
$ws.Cells.Item(61, 5).Value = "적"
$ws.Cells.Item(62, 5).Value = "체력이 가장 높은 아군"
$ws.Cells.Item(63, 5).Value = "체력이 가장 낮은 아군"
